# Scheduled runner update: refresh market-price derived columns (H-N) for
# several leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 812.6875
$ws.Range("I41").Value = 1432.2727
$ws.Range("J41").Value = 488.14285
$ws.Range("K41").Value = 1432.2727
$ws.Range("L41").Value = 488.14285
$ws.Range("M41").Value = -992.2727
$ws.Range("N41").Value = -1368.14285

$ws.Range("H106").Value = 3169.3845
$ws.Range("I106").Value = 4939.8
$ws.Range("J106").Value = 2062.875
$ws.Range("K106").Value = 4939.8
$ws.Range("L106").Value = 2062.875
$ws.Range("M106").Value = -4308.8
$ws.Range("N106").Value = -3324.875

$ws.Range("H129").Value = 2271.5625
$ws.Range("J129").Value = 1130.3846
$ws.Range("L129").Value = 3391.1538
$ws.Range("N129").Value = -13391.1538

$ws.Range("H137").Value = 982.4516
$ws.Range("I137").Value = 965.6429000000001
$ws.Range("K137").Value = 2896.9287
$ws.Range("M137").Value = -346.9287000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19813.734
$ws.Range("I32").Value = 4383.843
$ws.Range("J32").Value = 172398.22
$ws.Range("K32").Value = 4383.843
$ws.Range("L32").Value = 172398.22
$ws.Range("M32").Value = -4096.843
$ws.Range("N32").Value = -172972.22

$ws.Range("H61").Value = 1419.6666
$ws.Range("I61").Value = 1247.4348
$ws.Range("K61").Value = 1247.4348
$ws.Range("M61").Value = -1035.4348

$ws.Range("H74").Value = 720.1818
$ws.Range("I74").Value = 481.25
$ws.Range("J74").Value = 856.7143
$ws.Range("K74").Value = 481.25
$ws.Range("L74").Value = 856.7143
$ws.Range("M74").Value = 392.75
$ws.Range("N74").Value = -2604.7143

$ws.Range("H77").Value = 720.1818
$ws.Range("I77").Value = 481.25
$ws.Range("J77").Value = 856.7143
$ws.Range("K77").Value = 2406.25
$ws.Range("L77").Value = 4283.5715
$ws.Range("M77").Value = 1961.75
$ws.Range("N77").Value = -13019.5715

$ws.Range("H136").Value = 1419.6666
$ws.Range("I136").Value = 1247.4348
$ws.Range("K136").Value = 3742.3044
$ws.Range("M136").Value = -1192.3044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 2562.359
$ws.Range("I134").Value = 2318.1765
$ws.Range("J134").Value = 4222.8
$ws.Range("K134").Value = 6954.529500000001
$ws.Range("L134").Value = 12668.4
$ws.Range("M134").Value = -4419.529500000001
$ws.Range("N134").Value = -17738.4

$ws.Range("H135").Value = 52833.332
$ws.Range("J135").Value = 52833.332
$ws.Range("L135").Value = 52833.332
$ws.Range("N135").Value = -62973.332

$ws.Range("H137").Value = 39999
$ws.Range("J137").Value = 39999
$ws.Range("L137").Value = 39999
$ws.Range("N137").Value = -50199

$ws.Range("H140").Value = 65222.5
$ws.Range("J140").Value = 65222.5
$ws.Range("L140").Value = 65222.5
$ws.Range("N140").Value = -75582.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 42171.75
$ws.Range("I31").Value = 2297.2
$ws.Range("J31").Value = 70653.57000000001
$ws.Range("K31").Value = 2297.2
$ws.Range("L31").Value = 70653.57000000001
$ws.Range("M31").Value = -2002.2
$ws.Range("N31").Value = -71243.57000000001

$ws.Range("H34").Value = 42171.75
$ws.Range("I34").Value = 2297.2
$ws.Range("J34").Value = 70653.57000000001
$ws.Range("K34").Value = 2297.2
$ws.Range("L34").Value = 70653.57000000001
$ws.Range("M34").Value = -2095.2
$ws.Range("N34").Value = -71057.57000000001

$ws.Range("H134").Value = 1208.2963
$ws.Range("I134").Value = 1206.1666
$ws.Range("J134").Value = 1212.5555
$ws.Range("K134").Value = 3618.4998
$ws.Range("L134").Value = 3637.6665
$ws.Range("M134").Value = -1083.4998
$ws.Range("N134").Value = -8707.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1003.2439
$ws.Range("I5").Value = 563.4761999999999
$ws.Range("J5").Value = 1465
$ws.Range("K5").Value = 1690.4286
$ws.Range("L5").Value = 4395
$ws.Range("M5").Value = -1578.4286
$ws.Range("N5").Value = -4619

$ws.Range("H97").Value = 1099.75
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 1599.5
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 4798.5
$ws.Range("M97").Value = -1304
$ws.Range("N97").Value = -5790.5

$ws.Range("H111").Value = 2543.5
$ws.Range("I111").Value = 1999.3334
$ws.Range("J111").Value = 2870
$ws.Range("K111").Value = 5998.0002
$ws.Range("L111").Value = 8610
$ws.Range("M111").Value = -2931.0002
$ws.Range("N111").Value = -14744

$ws.Range("H112").Value = 68615.734
$ws.Range("I112").Value = 200845.4
$ws.Range("K112").Value = 602536.2
$ws.Range("M112").Value = -601428.2

$ws.Range("H131").Value = 1351.7283
$ws.Range("J131").Value = 1347.216
$ws.Range("L131").Value = 4041.648
$ws.Range("N131").Value = -14121.648

$ws.Range("H135").Value = 1003.2439
$ws.Range("I135").Value = 563.4761999999999
$ws.Range("J135").Value = 1465
$ws.Range("K135").Value = 5071.2858
$ws.Range("L135").Value = 13185
$ws.Range("M135").Value = -2536.2858
$ws.Range("N135").Value = -18255

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 44659.285
$ws.Range("J139").Value = 44659.285
$ws.Range("L139").Value = 44659.285
$ws.Range("N139").Value = -54939.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3890
$ws.Range("I100").Value = 2500
$ws.Range("J100").Value = 5280
$ws.Range("K100").Value = 2500
$ws.Range("L100").Value = 5280
$ws.Range("M100").Value = -1959
$ws.Range("N100").Value = -6362

$ws.Range("H136").Value = 982.01886
$ws.Range("I136").Value = 842.9400000000001
$ws.Range("J136").Value = 3300
$ws.Range("K136").Value = 2528.82
$ws.Range("L136").Value = 9900
$ws.Range("M136").Value = 21.17999999999984
$ws.Range("N136").Value = -15000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 43864.086
$ws.Range("I107").Value = 274.73334
$ws.Range("K107").Value = 824.20002
$ws.Range("M107").Value = 1095.79998

$ws.Range("H113").Value = 819.5833
$ws.Range("I113").Value = 796.5454999999999
$ws.Range("J113").Value = 839.0769
$ws.Range("K113").Value = 2389.6365
$ws.Range("L113").Value = 2517.2307
$ws.Range("M113").Value = -219.6364999999996
$ws.Range("N113").Value = -6857.2307

$ws.Range("H136").Value = 548.55
$ws.Range("I136").Value = 381.18182
$ws.Range("J136").Value = 1008.8125
$ws.Range("K136").Value = 1143.54546
$ws.Range("L136").Value = 3026.4375
$ws.Range("M136").Value = 1406.45454
$ws.Range("N136").Value = -8126.4375

Write-Host "Updated market-price columns on ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
